$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 5, 6) have had their species-record payload
# cyclically rotated: row2 -> row5, row5 -> row6, row6 -> row2.
# Capture the "before" values for the columns that carry the payload
# before overwriting anything.

$cols = @("A","B","E","F","G","H","Q","R","Y","AA")
$rows = @(2, 5, 6)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($col in $cols) {
        # NOTE: read via Value() (method-call syntax) rather than the bare
        # `.Value` property — the COM-interop shim's reflection resolves the
        # bare property read ambiguously here and hands back a descriptor
        # string instead of the cell's actual value.
        $snapshot[$r][$col] = $ws.Range("$col$r").Value()
    }
}

# Mapping of destination row -> source row (content moves FROM source TO destination)
$srcFor = @{ 2 = 6; 5 = 2; 6 = 5 }

foreach ($destRow in $rows) {
    $srcRow = $srcFor[$destRow]
    foreach ($col in $cols) {
        $value = $snapshot[$srcRow][$col]
        $cell = $ws.Range("$col$destRow")
        if ($col -eq "Y" -or $col -eq "AA") {
            # These columns store plain date-looking text (e.g. "2023-08-06"),
            # not real dates. Force text so Excel doesn't reinterpret it as
            # a date serial number.
            $cell.Value = "'" + $value
        } else {
            $cell.Value = $value
        }
    }
}
